# Auto-generated Excel COM-interop script applying the TPM update diff
# to the Cd14-Itgb1 LR-pairs sheet (NatmiData TPM pipeline output).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.6209125
$ws.Range("H2").Value = 3.241825
$ws.Range("I2").Value = 0.003190269591522575
$ws.Range("J2").Value = 0.002130042958301246
$ws.Range("M2").Value = 201.5557555
$ws.Range("N2").Value = 403.111511
$ws.Range("O2").Value = 0.2814680640969941
$ws.Range("P2").Value = 0.2271476997658372
$ws.Range("Q2").Value = 326.7042435368937
$ws.Range("R2").Value = 1306.816974147575
$ws.Range("S2").Value = 0.0008979590058733674
$ws.Range("T2").Value = 0.000483834358380547
# Row 3
$ws.Range("G3").Value = 1.6209125
$ws.Range("H3").Value = 3.241825
$ws.Range("I3").Value = 0.003190269591522575
$ws.Range("J3").Value = 0.002130042958301246
$ws.Range("N3").Value = 356.292984
$ws.Range("O3").Value = 0.1658517039268592
$ws.Range("P3").Value = 0.2007661144618275
$ws.Range("Q3").Value = 192.5065838093
$ws.Range("R3").Value = 1155.0395028558
$ws.Range("S3").Value = 0.0005291116477400642
$ws.Range("T3").Value = 0.0004276404483749176
# Row 4
$ws.Range("G4").Value = 1.6209125
$ws.Range("H4").Value = 3.241825
$ws.Range("I4").Value = 0.003190269591522575
$ws.Range("J4").Value = 0.002130042958301246
$ws.Range("M4").Value = 77.61293766666667
$ws.Range("N4").Value = 232.838813
$ws.Range("O4").Value = 0.1083847159795808
$ws.Range("P4").Value = 0.1312014153551619
$ws.Range("Q4").Value = 125.8037808256208
$ws.Range("R4").Value = 754.822684953725
$ws.Range("S4").Value = 0.0003457764635754677
$ws.Range("T4").Value = 0.0002794646508964195
# Row 5
$ws.Range("G5").Value = 1.6209125
$ws.Range("H5").Value = 3.241825
$ws.Range("I5").Value = 0.003190269591522575
$ws.Range("J5").Value = 0.002130042958301246
$ws.Range("M5").Value = 172.0397415
$ws.Range("N5").Value = 344.079483
$ws.Range("O5").Value = 0.2402496166265631
$ws.Range("P5").Value = 0.1938839774289365
$ws.Range("Q5").Value = 278.8613674941187
$ws.Range("R5").Value = 1115.445469976475
$ws.Range("S5").Value = 0.0007664610462986808
$ws.Range("T5").Value = 0.0004129812008499438
# Row 6
$ws.Range("G6").Value = 1.6209125
$ws.Range("H6").Value = 3.241825
$ws.Range("I6").Value = 0.003190269591522575
$ws.Range("J6").Value = 0.002130042958301246
$ws.Range("M6").Value = 59.85226566666666
$ws.Range("N6").Value = 179.556797
$ws.Range("O6").Value = 0.08358233833226186
$ws.Range("P6").Value = 0.1011777443782085
$ws.Range("Q6").Value = 97.01528557242082
$ws.Range("R6").Value = 582.091713434525
$ws.Range("S6").Value = 0.0002666501923697667
$ws.Range("T6").Value = 0.0002155129419496064
# Row 7
$ws.Range("G7").Value = 1.6209125
$ws.Range("H7").Value = 3.241825
$ws.Range("I7").Value = 0.003190269591522575
$ws.Range("J7").Value = 0.002130042958301246
$ws.Range("M7").Value = 86.26244733333333
$ws.Range("N7").Value = 258.787342
$ws.Range("O7").Value = 0.1204635610377409
$ws.Range("P7").Value = 0.1458230486100285
$ws.Range("Q7").Value = 139.8238791631917
$ws.Range("R7").Value = 838.9432749791499
$ws.Range("S7").Value = 0.0003843112356652286
$ws.Range("T7").Value = 0.0003106093578498114
# Row 8
$ws.Range("I8").Value = 0.001893391735996584
$ws.Range("J8").Value = 0.001896237427057309
$ws.Range("M8").Value = 201.5557555
$ws.Range("N8").Value = 403.111511
$ws.Range("O8").Value = 0.2814680640969941
$ws.Range("P8").Value = 0.2271476997658372
$ws.Range("Q8").Value = 193.8955618269707
$ws.Range("R8").Value = 1163.373370961824
$ws.Range("S8").Value = 0.0005329293065082053
$ws.Range("T8").Value = 0.0004307259697659574
# Row 9
$ws.Range("I9").Value = 0.001893391735996584
$ws.Range("J9").Value = 0.001896237427057309
$ws.Range("N9").Value = 356.292984
$ws.Range("O9").Value = 0.1658517039268592
$ws.Range("P9").Value = 0.2007661144618275
$ws.Range("Q9").Value = 114.2506501262507
$ws.Range("S9").Value = 0.0003140222456160674
$ws.Range("T9").Value = 0.0003807002203273891
# Row 10
$ws.Range("I10").Value = 0.001893391735996584
$ws.Range("J10").Value = 0.001896237427057309
$ws.Range("M10").Value = 77.61293766666667
$ws.Range("N10").Value = 232.838813
$ws.Range("O10").Value = 0.1083847159795808
$ws.Range("P10").Value = 0.1312014153551619
$ws.Range("Q10").Value = 74.66323209966578
$ws.Range("R10").Value = 671.969088896992
$ws.Range("S10").Value = 0.0002052147255440752
$ws.Range("T10").Value = 0.0002487890342793496
# Row 11
$ws.Range("I11").Value = 0.001893391735996584
$ws.Range("J11").Value = 0.001896237427057309
$ws.Range("M11").Value = 172.0397415
$ws.Range("N11").Value = 344.079483
$ws.Range("O11").Value = 0.2402496166265631
$ws.Range("P11").Value = 0.1938839774289365
$ws.Range("Q11").Value = 165.501313777712
$ws.Range("R11").Value = 993.0078826662719
$ws.Range("S11").Value = 0.000454886638697082
$ws.Range("T11").Value = 0.000367650054507484
# Row 12
$ws.Range("I12").Value = 0.001893391735996584
$ws.Range("J12").Value = 0.001896237427057309
$ws.Range("M12").Value = 59.85226566666666
$ws.Range("N12").Value = 179.556797
$ws.Range("O12").Value = 0.08358233833226186
$ws.Range("P12").Value = 0.1011777443782085
$ws.Range("Q12").Value = 57.57756035924977
$ws.Range("R12").Value = 518.198043233248
$ws.Range("S12").Value = 0.0001582541086735751
$ws.Range("T12").Value = 0.0001918570256751962
# Row 13
$ws.Range("I13").Value = 0.001893391735996584
$ws.Range("J13").Value = 0.001896237427057309
$ws.Range("M13").Value = 86.26244733333333
$ws.Range("N13").Value = 258.787342
$ws.Range("O13").Value = 0.1204635610377409
$ws.Range("P13").Value = 0.1458230486100285
$ws.Range("Q13").Value = 82.98401426828089
$ws.Range("R13").Value = 746.8561284145279
$ws.Range("S13").Value = 0.0002280847109575788
$ws.Range("T13").Value = 0.0002765151225019333
# Row 14
$ws.Range("G14").Value = 209.9667966666667
$ws.Range("H14").Value = 629.90039
$ws.Range("I14").Value = 0.4132553031919183
$ws.Range("J14").Value = 0.4138764091678941
$ws.Range("M14").Value = 201.5557555
$ws.Range("N14").Value = 403.111511
$ws.Range("O14").Value = 0.2814680640969941
$ws.Range("P14").Value = 0.2271476997658372
$ws.Range("Q14").Value = 42320.01633206489
$ws.Range("R14").Value = 253920.0979923893
$ws.Range("S14").Value = 0.1163181701672456
$ws.Range("T14").Value = 0.0940110743298316
# Row 15
$ws.Range("G15").Value = 209.9667966666667
$ws.Range("H15").Value = 629.90039
$ws.Range("I15").Value = 0.4132553031919183
$ws.Range("J15").Value = 0.4138764091678941
$ws.Range("N15").Value = 356.292984
$ws.Range("O15").Value = 0.1658517039268592
$ws.Range("P15").Value = 0.2007661144618275
$ws.Range("Q15").Value = 24936.56550842931
$ws.Range("R15").Value = 224429.0895758638
$ws.Range("S15").Value = 0.06853909619119047
$ws.Range("T15").Value = 0.0830923585360516
# Row 16
$ws.Range("G16").Value = 209.9667966666667
$ws.Range("H16").Value = 629.90039
$ws.Range("I16").Value = 0.4132553031919183
$ws.Range("J16").Value = 0.4138764091678941
$ws.Range("M16").Value = 77.61293766666667
$ws.Range("N16").Value = 232.838813
$ws.Range("O16").Value = 0.1083847159795808
$ws.Range("P16").Value = 0.1312014153551619
$ws.Range("Q16").Value = 16296.13990175968
$ws.Range("R16").Value = 146665.2591158371
$ws.Range("S16").Value = 0.04479055866351164
$ws.Range("T16").Value = 0.05430117066493982
# Row 17
$ws.Range("G17").Value = 209.9667966666667
$ws.Range("H17").Value = 629.90039
$ws.Range("I17").Value = 0.4132553031919183
$ws.Range("J17").Value = 0.4138764091678941
$ws.Range("M17").Value = 172.0397415
$ws.Range("N17").Value = 344.079483
$ws.Range("O17").Value = 0.2402496166265631
$ws.Range("P17").Value = 0.1938839774289365
$ws.Range("Q17").Value = 36122.63342211639
$ws.Range("R17").Value = 216735.8005326984
$ws.Range("S17").Value = 0.09928442816075247
$ws.Range("T17").Value = 0.08024400437347727
# Row 18
$ws.Range("G18").Value = 209.9667966666667
$ws.Range("H18").Value = 629.90039
$ws.Range("I18").Value = 0.4132553031919183
$ws.Range("J18").Value = 0.4138764091678941
$ws.Range("M18").Value = 59.85226566666666
$ws.Range("N18").Value = 179.556797
$ws.Range("O18").Value = 0.08358233833226186
$ws.Range("P18").Value = 0.1011777443782085
$ws.Range("Q18").Value = 12566.98849527231
$ws.Range("R18").Value = 113102.8964574508
$ws.Range("S18").Value = 0.03454084456898837
$ws.Range("T18").Value = 0.04187508153096001
# Row 19
$ws.Range("G19").Value = 209.9667966666667
$ws.Range("H19").Value = 629.90039
$ws.Range("I19").Value = 0.4132553031919183
$ws.Range("J19").Value = 0.4138764091678941
$ws.Range("M19").Value = 86.26244733333333
$ws.Range("N19").Value = 258.787342
$ws.Range("O19").Value = 0.1204635610377409
$ws.Range("P19").Value = 0.1458230486100285
$ws.Range("Q19").Value = 18112.24973920704
$ws.Range("R19").Value = 163010.2476528634
$ws.Range("S19").Value = 0.04978220544022979
$ws.Range("T19").Value = 0.06035271973263386
# Row 20
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.6665209999999999
$ws.Range("H20").Value = 1.333042
$ws.Range("I20").Value = 0.001311842359418672
$ws.Range("J20").Value = 0.000875876003553495
$ws.Range("M20").Value = 201.5557555
$ws.Range("N20").Value = 403.111511
$ws.Range("O20").Value = 0.2814680640969941
$ws.Range("P20").Value = 0.2271476997658372
$ws.Range("Q20").Value = 134.3411437116155
$ws.Range("R20").Value = 537.3645748464619
$ws.Range("S20").Value = 0.0003692417293060067
$ws.Range("T20").Value = 0.0001989532194872706
# Row 21
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.6665209999999999
$ws.Range("H21").Value = 1.333042
$ws.Range("I21").Value = 0.001311842359418672
$ws.Range("J21").Value = 0.000875876003553495
$ws.Range("N21").Value = 356.292984
$ws.Range("O21").Value = 0.1658517039268592
$ws.Range("P21").Value = 0.2007661144618275
$ws.Range("Q21").Value = 79.15891866288798
$ws.Range("R21").Value = 474.9535119773279
$ws.Range("S21").Value = 0.000217571290593018
$ws.Range("T21").Value = 0.000175846221983789
# Row 22
$ws.Range("E22").Value = 2
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 0.6665209999999999
$ws.Range("H22").Value = 1.333042
$ws.Range("I22").Value = 0.001311842359418672
$ws.Range("J22").Value = 0.000875876003553495
$ws.Range("M22").Value = 77.61293766666667
$ws.Range("N22").Value = 232.838813
$ws.Range("O22").Value = 0.1083847159795808
$ws.Range("P22").Value = 0.1312014153551619
$ws.Range("Q22").Value = 51.73065282652433
$ws.Range("R22").Value = 310.383916959146
$ws.Range("S22").Value = 0.000142183661535576
$ws.Range("T22").Value = 0.0001149161713418414
# Row 23
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 0.6665209999999999
$ws.Range("H23").Value = 1.333042
$ws.Range("I23").Value = 0.001311842359418672
$ws.Range("J23").Value = 0.000875876003553495
$ws.Range("M23").Value = 172.0397415
$ws.Range("N23").Value = 344.079483
$ws.Range("O23").Value = 0.2402496166265631
$ws.Range("P23").Value = 0.1938839774289365
$ws.Range("Q23").Value = 114.6681005443215
$ws.Range("R23").Value = 458.6724021772859
$ws.Range("S23").Value = 0.000315169623924822
$ws.Range("T23").Value = 0.0001698183233035129
# Row 24
$ws.Range("E24").Value = 2
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 0.6665209999999999
$ws.Range("H24").Value = 1.333042
$ws.Range("I24").Value = 0.001311842359418672
$ws.Range("J24").Value = 0.000875876003553495
$ws.Range("M24").Value = 59.85226566666666
$ws.Range("N24").Value = 179.556797
$ws.Range("O24").Value = 0.08358233833226186
$ws.Range("P24").Value = 0.1011777443782085
$ws.Range("Q24").Value = 39.89279196441233
$ws.Range("R24").Value = 239.356751786474
$ws.Range("S24").Value = 0.0001096468519235241
$ws.Range("T24").Value = 0.00008861915839454234
# Row 25
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 0.6665209999999999
$ws.Range("H25").Value = 1.333042
$ws.Range("I25").Value = 0.001311842359418672
$ws.Range("J25").Value = 0.000875876003553495
$ws.Range("M25").Value = 86.26244733333333
$ws.Range("N25").Value = 258.787342
$ws.Range("O25").Value = 0.1204635610377409
$ws.Range("P25").Value = 0.1458230486100285
$ws.Range("Q25").Value = 57.49573265906066
$ws.Range("R25").Value = 344.9743959543639
$ws.Range("S25").Value = 0.0001580292021357253
$ws.Range("T25").Value = 0.0001277229090425388
# Row 26
$ws.Range("E26").Value = 3
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 169.4444526666666
$ws.Range("H26").Value = 508.333358
$ws.Range("I26").Value = 0.3334994855025505
$ws.Range("J26").Value = 0.334000721716171
$ws.Range("M26").Value = 201.5557555
$ws.Range("N26").Value = 403.111511
$ws.Range("O26").Value = 0.2814680640969941
$ws.Range("P26").Value = 0.2271476997658372
$ws.Range("Q26").Value = 34152.50467251398
$ws.Range("R26").Value = 204915.0280350839
$ws.Range("S26").Value = 0.09386945456174645
$ws.Range("T26").Value = 0.07586749565795775
# Row 27
$ws.Range("E27").Value = 3
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 169.4444526666666
$ws.Range("H27").Value = 508.333358
$ws.Range("I27").Value = 0.3334994855025505
$ws.Range("J27").Value = 0.334000721716171
$ws.Range("N27").Value = 356.292984
$ws.Range("O27").Value = 0.1658517039268592
$ws.Range("P27").Value = 0.2007661144618275
$ws.Range("Q27").Value = 20123.95655428447
$ws.Range("R27").Value = 181115.6089885603
$ws.Range("S27").Value = 0.05531145792932889
$ws.Range("T27").Value = 0.06705602712640178
# Row 28
$ws.Range("E28").Value = 3
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 169.4444526666666
$ws.Range("H28").Value = 508.333358
$ws.Range("I28").Value = 0.3334994855025505
$ws.Range("J28").Value = 0.334000721716171
$ws.Range("M28").Value = 77.61293766666667
$ws.Range("N28").Value = 232.838813
$ws.Range("O28").Value = 0.1083847159795808
$ws.Range("P28").Value = 0.1312014153551619
$ws.Range("Q28").Value = 13151.08174278045
$ws.Range("R28").Value = 118359.7356850241
$ws.Range("S28").Value = 0.03614624701553028
$ws.Range("T28").Value = 0.0438213674188072
# Row 29
$ws.Range("E29").Value = 3
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 169.4444526666666
$ws.Range("H29").Value = 508.333358
$ws.Range("I29").Value = 0.3334994855025505
$ws.Range("J29").Value = 0.334000721716171
$ws.Range("M29").Value = 172.0397415
$ws.Range("N29").Value = 344.079483
$ws.Range("O29").Value = 0.2402496166265631
$ws.Range("P29").Value = 0.1938839774289365
$ws.Range("Q29").Value = 29151.17983538231
$ws.Range("R29").Value = 174907.0790122939
$ws.Range("S29").Value = 0.0801231235371438
$ws.Range("T29").Value = 0.0647573883904666
# Row 30
$ws.Range("E30").Value = 3
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 169.4444526666666
$ws.Range("H30").Value = 508.333358
$ws.Range("I30").Value = 0.3334994855025505
$ws.Range("J30").Value = 0.334000721716171
$ws.Range("M30").Value = 59.85226566666666
$ws.Range("N30").Value = 179.556797
$ws.Range("O30").Value = 0.08358233833226186
$ws.Range("P30").Value = 0.1011777443782085
$ws.Range("Q30").Value = 10141.63439674826
$ws.Range("R30").Value = 91274.70957073431
$ws.Range("S30").Value = 0.02787466683090944
$ws.Range("T30").Value = 0.0337934396439359
# Row 31
$ws.Range("E31").Value = 3
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 169.4444526666666
$ws.Range("H31").Value = 508.333358
$ws.Range("I31").Value = 0.3334994855025505
$ws.Range("J31").Value = 0.334000721716171
$ws.Range("M31").Value = 86.26244733333333
$ws.Range("N31").Value = 258.787342
$ws.Range("O31").Value = 0.1204635610377409
$ws.Range("P31").Value = 0.1458230486100285
$ws.Range("Q31").Value = 14616.69317408382
$ws.Range("R31").Value = 131550.2385667544
$ws.Range("S31").Value = 0.0401745356278917
$ws.Range("T31").Value = 0.04870500347860179
# Row 32
$ws.Range("G32").Value = 125.4194246666667
$ws.Range("H32").Value = 376.258274
$ws.Range("I32").Value = 0.2468497076185933
$ws.Range("J32").Value = 0.2472207127270228
$ws.Range("M32").Value = 201.5557555
$ws.Range("N32").Value = 403.111511
$ws.Range("O32").Value = 0.2814680640969941
$ws.Range("P32").Value = 0.2271476997658372
$ws.Range("Q32").Value = 25279.00689306533
$ws.Range("R32").Value = 151674.041358392
$ws.Range("S32").Value = 0.06948030932631447
$ws.Range("T32").Value = 0.05615561623041406
# Row 33
$ws.Range("G33").Value = 125.4194246666667
$ws.Range("H33").Value = 376.258274
$ws.Range("I33").Value = 0.2468497076185933
$ws.Range("J33").Value = 0.2472207127270228
$ws.Range("N33").Value = 356.292984
$ws.Range("O33").Value = 0.1658517039268592
$ws.Range("P33").Value = 0.2007661144618275
$ws.Range("Q33").Value = 14895.35368868329
$ws.Range("R33").Value = 134058.1831981496
$ws.Range("S33").Value = 0.04094044462239069
$ws.Range("T33").Value = 0.04963354190868804
# Row 34
$ws.Range("G34").Value = 125.4194246666667
$ws.Range("H34").Value = 376.258274
$ws.Range("I34").Value = 0.2468497076185933
$ws.Range("J34").Value = 0.2472207127270228
$ws.Range("M34").Value = 77.61293766666667
$ws.Range("N34").Value = 232.838813
$ws.Range("O34").Value = 0.1083847159795808
$ws.Range("P34").Value = 0.1312014153551619
$ws.Range("Q34").Value = 9734.169988843196
$ws.Range("R34").Value = 87607.52989958876
$ws.Range("S34").Value = 0.0267547354498838
$ws.Range("T34").Value = 0.03243570741489728
# Row 35
$ws.Range("G35").Value = 125.4194246666667
$ws.Range("H35").Value = 376.258274
$ws.Range("I35").Value = 0.2468497076185933
$ws.Range("J35").Value = 0.2472207127270228
$ws.Range("M35").Value = 172.0397415
$ws.Range("N35").Value = 344.079483
$ws.Range("O35").Value = 0.2402496166265631
$ws.Range("P35").Value = 0.1938839774289365
$ws.Range("Q35").Value = 21577.12539873206
$ws.Range("R35").Value = 129462.7523923923
$ws.Range("S35").Value = 0.05930554761974623
$ws.Range("T35").Value = 0.04793213508633168
# Row 36
$ws.Range("G36").Value = 125.4194246666667
$ws.Range("H36").Value = 376.258274
$ws.Range("I36").Value = 0.2468497076185933
$ws.Range("J36").Value = 0.2472207127270228
$ws.Range("M36").Value = 59.85226566666666
$ws.Range("N36").Value = 179.556797
$ws.Range("O36").Value = 0.08358233833226186
$ws.Range("P36").Value = 0.1011777443782085
$ws.Range("Q36").Value = 7506.636724909818
$ws.Range("R36").Value = 67559.73052418837
$ws.Range("S36").Value = 0.02063227577939718
$ws.Range("T36").Value = 0.02501323407729322
# Row 37
$ws.Range("G37").Value = 125.4194246666667
$ws.Range("H37").Value = 376.258274
$ws.Range("I37").Value = 0.2468497076185933
$ws.Range("J37").Value = 0.2472207127270228
$ws.Range("M37").Value = 86.26244733333333
$ws.Range("N37").Value = 258.787342
$ws.Range("O37").Value = 0.1204635610377409
$ws.Range("P37").Value = 0.1458230486100285
$ws.Range("Q37").Value = 10818.9865148853
$ws.Range("R37").Value = 97370.87863396769
$ws.Range("S37").Value = 0.02973639482086092
$ws.Range("T37").Value = 0.03605047800939852
